$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 159, shifting existing rows 159:264 down to 160:265
$ws.Rows.Item(159).Insert()

# Populate the newly inserted row 159 with the new weekly record
$ws.Cells.Item(159, 1).Value = 11
$ws.Cells.Item(159, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(159, 3).Value = "Bíobío"
$ws.Cells.Item(159, 4).Value = 44603
$ws.Cells.Item(159, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(159, 5).Value = 8
$ws.Cells.Item(159, 6).Value = 100114014
$ws.Cells.Item(159, 7).Value = "Betarraga"
$ws.Cells.Item(159, 8).Value = "Sin especificar"
$ws.Cells.Item(159, 9).Value = "Primera"
$ws.Cells.Item(159, 10).Value = 1600
$ws.Cells.Item(159, 11).Value = 600
$ws.Cells.Item(159, 12).Value = 650
$ws.Cells.Item(159, 13).Value = 625
$ws.Cells.Item(159, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(159, 15).Value = "Región Metropolitana"
$ws.Cells.Item(159, 16).Value = 125
$ws.Cells.Item(159, 17).Value = 5
$ws.Cells.Item(159, 18).Value = "Hortaliza"
